$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.47"
$ws.Range("E2").Value = "'2.55%"
$ws.Range("D3").Value = "'35.03"
$ws.Range("E3").Value = "'12.61%"
$ws.Range("D4").Value = "'5.132"
$ws.Range("E4").Value = "'4.34%"
$ws.Range("D5").Value = "'0.07767"
$ws.Range("E5").Value = "'4.39%"
$ws.Range("D6").Value = "'2.345"
$ws.Range("E6").Value = "'5.65%"
$ws.Range("D7").Value = "'8.036"
$ws.Range("E7").Value = "'3.71%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'3.944"
$ws.Range("E8").Value = "'5.37%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9312"
$ws.Range("E9").Value = "'1.58%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1000"
$ws.Range("E10").Value = "'11.82%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1791"
$ws.Range("E11").Value = "'4.53%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08623"
$ws.Range("E12").Value = "'3.88%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03322"
$ws.Range("E13").Value = "'6.66%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09923"
$ws.Range("E14").Value = "'-1.56%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001500"
$ws.Range("E15").Value = "'-1.14%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005776"
$ws.Range("E16").Value = "'0.28%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.460"
$ws.Range("E17").Value = "'-1.27%"
$ws.Range("D18").Value = "'2.139"
$ws.Range("E18").Value = "'3.06%"
$ws.Range("D19").Value = "'0.3366"
$ws.Range("E19").Value = "'1.11%"
$ws.Range("E20").Value = "'2.75%"
$ws.Range("D21").Value = "'4.303"
$ws.Range("E21").Value = "'8.04%"
$ws.Range("E23").Value = "'-0.56%"
$ws.Range("E24").Value = "'0.25%"
$ws.Range("D25").Value = "'0.004378"
$ws.Range("E25").Value = "'-5.20%"
$ws.Range("E26").Value = "'-0.02%"
$ws.Range("E27").Value = "'-0.09%"
$ws.Range("D39").Value = "'0.01793"
$ws.Range("E39").Value = "'11.28%"
$ws.Range("D40").Value = "'0.04797"
$ws.Range("E40").Value = "'6.95%"
$ws.Range("D41").Value = "'0.007784"
$ws.Range("E41").Value = "'6.81%"
$ws.Range("D42").Value = "'0.1413"
$ws.Range("E42").Value = "'6.20%"
$ws.Range("D43").Value = "'0.007227"
$ws.Range("E43").Value = "'-19.33%"
$ws.Range("D44").Value = "'0.002073"
$ws.Range("E44").Value = "'5.50%"
$ws.Range("D45").Value = "'0.009456"
$ws.Range("E45").Value = "'9.79%"
$ws.Range("D46").Value = "'0.00006112"
$ws.Range("E46").Value = "'0.53%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.02%"
$ws.Range("D48").Value = "'3.027"
$ws.Range("E48").Value = "'35.69%"
$ws.Range("E49").Value = "'-0.03%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.02%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.02%"
